$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shift rows 24..30 down to 25..31 (bottom-up so we never overwrite a row
#    before we've read it). We rebuild each cell explicitly (value/formula +
#    style) rather than using Range.Insert()/Range.Copy() because both of
#    those blow the shared HYPERLINK formulas in columns D/E out into
#    per-cell formulas and/or duplicate the row across all 16384 columns.
# ---------------------------------------------------------------------------

for ($r = 30; $r -ge 24; $r--) {
    $dst = $r + 1

    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $f = $ws.Cells.Item($r, 6).Value2
    $g = $ws.Cells.Item($r, 7).Value2
    $h = $ws.Cells.Item($r, 8).Value2
    $i = $ws.Cells.Item($r, 9).Value2
    $j = $ws.Cells.Item($r, 10).Value2
    $k = $ws.Cells.Item($r, 11).Value2

    $hasD = $ws.Cells.Item($r, 4).HasFormula
    $hasE = $ws.Cells.Item($r, 5).HasFormula

    if ($null -ne $a) { $ws.Cells.Item($dst, 1).Value = $a } else { $ws.Cells.Item($dst, 1).ClearContents() }
    if ($null -ne $b) { $ws.Cells.Item($dst, 2).Value = $b } else { $ws.Cells.Item($dst, 2).ClearContents() }
    if ($null -ne $c) { $ws.Cells.Item($dst, 3).Value = $c } else { $ws.Cells.Item($dst, 3).ClearContents() }

    if ($hasD) {
        $ws.Cells.Item($dst, 4).Formula = "=HYPERLINK(K$dst)"
    } elseif ($null -ne $ws.Cells.Item($r, 4).Value2) {
        $ws.Cells.Item($dst, 4).Value = $ws.Cells.Item($r, 4).Value2
    } else {
        $ws.Cells.Item($dst, 4).ClearContents()
    }

    if ($hasE) {
        $ws.Cells.Item($dst, 5).Formula = "=HYPERLINK(J$dst)"
    } elseif ($null -ne $ws.Cells.Item($r, 5).Value2) {
        $ws.Cells.Item($dst, 5).Value = $ws.Cells.Item($r, 5).Value2
    } else {
        $ws.Cells.Item($dst, 5).ClearContents()
    }

    if ($null -ne $f) { $ws.Cells.Item($dst, 6).Value = $f } else { $ws.Cells.Item($dst, 6).ClearContents() }
    if ($null -ne $g) { $ws.Cells.Item($dst, 7).Value = $g } else { $ws.Cells.Item($dst, 7).ClearContents() }
    if ($null -ne $h) { $ws.Cells.Item($dst, 8).Value = $h } else { $ws.Cells.Item($dst, 8).ClearContents() }
    if ($null -ne $i) { $ws.Cells.Item($dst, 9).Value = $i } else { $ws.Cells.Item($dst, 9).ClearContents() }
    if ($null -ne $j) { $ws.Cells.Item($dst, 10).Value = $j } else { $ws.Cells.Item($dst, 10).ClearContents() }
    if ($null -ne $k) { $ws.Cells.Item($dst, 11).Value = $k } else { $ws.Cells.Item($dst, 11).ClearContents() }
}

# ---------------------------------------------------------------------------
# 2. Write the new row 24 (New Vista / BVSD) that was inserted by the edit.
#    D24/E24 are plain literal strings in the target (not HYPERLINK formulas).
# ---------------------------------------------------------------------------

$ws.Cells.Item(24, 1).Value = "New Vista"
$ws.Cells.Item(24, 2).Value = "700 20th St, Boulder, CO 80302"
$ws.Cells.Item(24, 3).Value = "BVSD"
$ws.Cells.Item(24, 4).Value = "http://nvh.bvsd.org/"
$ws.Cells.Item(24, 5).Value = "https://maps.app.goo.gl/mVZ7Cv3jx5iSPW3u9"
$ws.Cells.Item(24, 6).Value = "New Vista High School"
$ws.Cells.Item(24, 7).Value = "POINT (-105.2665043792843 40.00119784626813)"
$ws.Cells.Item(24, 8).Value = 40.001197846268099
$ws.Cells.Item(24, 9).Value = -105.266504379284
$ws.Cells.Item(24, 10).ClearContents()
$ws.Cells.Item(24, 11).ClearContents()

# ---------------------------------------------------------------------------
# 3. Selection / view bookkeeping to mirror the saved UI state: the user's
#    last action selected the newly inserted row 24 (whole-row selection),
#    and scrolled the sheet back so column A is visible again.
# ---------------------------------------------------------------------------

$ws.Range("A24:XFD24").Select()
$ws.Activate()

Write-Output "done"
